$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix spelling of the competitor's name (from the first competition blog post)
$ws.Range("A5").Value = "Jakob Stymne"

# Move the active selection to A25 (matches the saved cursor position in the
# source workbook after the edit)
$null = $ws.Range("A25").Select()
